# Append 18 new transaction rows (rows 173-190) for the week ending 2021-02-21
# to the "Konto" worksheet of OrderReceipts_Expenses_Konto_Classification.xlsx

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Grab the date number format used by the existing date column (column A)
# so new date cells match the existing style exactly.
$dateFormat = $ws.Cells.Item(172, 1).NumberFormat()

# Each entry: Row, Date(serial), ReceiptNumber, Konto, Beskrivning, Debet, Kredit
$rows = @(
    @{ R=173; A=44242; B=$null;    C=4010; D="MATVARLDEN VEDD K0135";            E=156.88;            F=$null },
    @{ R=174; A=44242; B=$null;    C=2645; D="MATVARLDEN VEDD K0135";            E=18.82;             F=$null },
    @{ R=175; A=44242; B=$null;    C=1930; D="MATVARLDEN VEDD K0135";            E=$null;             F=175.7 },

    @{ R=176; A=44245; B=2181305;  C=3011; D="Order 2181305 Swish +46734333950"; E=$null;             F=691.0700000000001 },
    @{ R=177; A=44245; B=2181305;  C=2611; D="Order 2181305 Swish +46734333950"; E=$null;             F=82.93000000000001 },
    @{ R=178; A=44245; B=2181305;  C=1930; D="Order 2181305 Swish +46734333950"; E=774;                F=$null },

    @{ R=179; A=44245; B=4181944;  C=3011; D="Order 4181944 Swish +46763141239"; E=$null;             F=725.89 },
    @{ R=180; A=44245; B=4181944;  C=2611; D="Order 4181944 Swish +46763141239"; E=$null;             F=87.11 },
    @{ R=181; A=44245; B=4181944;  C=1930; D="Order 4181944 Swish +46763141239"; E=813;                F=$null },

    @{ R=182; A=44247; B=$null;    C=6400; D="FACEBK J3NSNYWY62 K6885";          E=415;               F=$null },
    @{ R=183; A=44247; B=$null;    C=$null;D="FACEBK J3NSNYWY62 K6885";          E=0;                 F=$null },
    @{ R=184; A=44247; B=$null;    C=6400; D="FACEBK YCYG9YSZ62 K6885";          E=3;                 F=$null },
    @{ R=185; A=44247; B=$null;    C=$null;D="FACEBK YCYG9YSZ62 K6885";          E=0;                 F=$null },
    @{ R=186; A=44247; B=$null;    C=1930; D="FACEBK YCYG9YSZ62 K6885";          E=$null;             F=3 },
    @{ R=187; A=44247; B=$null;    C=1930; D="FACEBK J3NSNYWY62 K6885";          E=$null;             F=415 },

    @{ R=188; A=44247; B=$null;    C=4010; D="NGROCERIES AB K0135";              E=204.29;            F=$null },
    @{ R=189; A=44247; B=$null;    C=2645; D="NGROCERIES AB K0135";              E=24.51;             F=$null },
    @{ R=190; A=44247; B=$null;    C=1930; D="NGROCERIES AB K0135";              E=$null;             F=228.8 }
)

foreach ($row in $rows) {
    $r = $row.R

    $cellA = $ws.Cells.Item($r, 1)
    $cellA.Value = $row.A
    $cellA.NumberFormat = $dateFormat

    $cellB = $ws.Cells.Item($r, 2)
    if ($null -eq $row.B) { $cellB.Value = "" } else { $cellB.Value = $row.B }

    $cellC = $ws.Cells.Item($r, 3)
    if ($null -eq $row.C) { $cellC.Value = "" } else { $cellC.Value = $row.C }

    $cellD = $ws.Cells.Item($r, 4)
    $cellD.Value = $row.D

    $cellE = $ws.Cells.Item($r, 5)
    if ($null -eq $row.E) { $cellE.Value = "" } else { $cellE.Value = $row.E }

    $cellF = $ws.Cells.Item($r, 6)
    if ($null -eq $row.F) { $cellF.Value = "" } else { $cellF.Value = $row.F }
}

Write-Host ("Added {0} rows; new dimension: {1}" -f $rows.Count, $ws.UsedRange.Address())
